$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.845013590853341
$ws.Range("C2").Value = 0.1537331888545168
$ws.Range("D2").Value = 0.07770624144853855
$ws.Range("E2").Value = 0.07113721758956615
$ws.Range("G2").Value = 2.124212706615339
$ws.Range("H2").Value = 1.671350690205969
$ws.Range("L2").Value = 0.2405452966810202
$ws.Range("B3").Value = 1.734008676065002
$ws.Range("C3").Value = 0.1359928575599554
$ws.Range("D3").Value = 0.07057916706918377
$ws.Range("E3").Value = 0.07122251944040059
$ws.Range("G3").Value = 2.054630021345929
$ws.Range("H3").Value = 1.645254660415617
$ws.Range("L3").Value = 0.2316589118509569
$ws.Range("B4").Value = 1.667079613234932
$ws.Range("C4").Value = 0.1250446443515045
$ws.Range("D4").Value = 0.06624682934017301
$ws.Range("E4").Value = 0.07129514232904555
$ws.Range("G4").Value = 2.013168216694027
$ws.Range("H4").Value = 1.630096018852299
$ws.Range("L4").Value = 0.2263594676026628
$ws.Range("B5").Value = 1.640112957842121
$ws.Range("C5").Value = 0.1205689255971265
$ws.Range("D5").Value = 0.06449215524688157
$ws.Range("E5").Value = 0.07132981841615127
$ws.Range("G5").Value = 1.996585776842409
$ws.Range("H5").Value = 1.624134465178088
$ws.Range("L5").Value = 0.2242391444272016
$ws.Range("B6").Value = 1.635653696407303
$ws.Range("C6").Value = 0.1198248689121613
$ws.Range("D6").Value = 0.06420143865980776
$ws.Range("E6").Value = 0.07133588299443083
$ws.Range("G6").Value = 1.993851103675382
$ws.Range("H6").Value = 1.623157534371558
$ws.Range("L6").Value = 0.2238894310350048
$ws.Range("B7").Value = 1.66671468834619
$ws.Range("C7").Value = 0.1249843412018947
$ws.Range("D7").Value = 0.06622312177329093
$ws.Range("E7").Value = 0.07129558941975844
$ws.Range("G7").Value = 2.012943315475951
$ws.Range("H7").Value = 1.630014748039088
$ws.Range("L7").Value = 0.2263307135373935
$ws.Range("B8").Value = 1.806483235213477
$ws.Range("C8").Value = 0.1476276957644131
$ws.Range("D8").Value = 0.07523962771844594
$ws.Range("E8").Value = 0.0711624214521116
$ws.Range("G8").Value = 2.099956289286638
$ws.Range("H8").Value = 1.662172359793061
$ws.Range("L8").Value = 0.2374485996890598
$ws.Range("B9").Value = 2.090398638235399
$ws.Range("C9").Value = 0.1916042013069443
$ws.Range("D9").Value = 0.09327800897877125
$ws.Range("E9").Value = 0.0710624056574396
$ws.Range("G9").Value = 2.280786537497107
$ws.Range("H9").Value = 1.732171573508879
$ws.Range("L9").Value = 0.2605054732274965
$ws.Range("B10").Value = 2.305126325548088
$ws.Range("C10").Value = 0.2236760781646296
$ws.Range("D10").Value = 0.1067640070761939
$ws.Range("E10").Value = 0.07108786264558553
$ws.Range("G10").Value = 2.420127632426215
$ws.Range("H10").Value = 1.787946806656777
$ws.Range("L10").Value = 0.2782269455129409
$ws.Range("B11").Value = 2.404177183460945
$ws.Range("C11").Value = 0.238220024434753
$ws.Range("D11").Value = 0.1129532542650367
$ws.Range("E11").Value = 0.07112108675379325
$ws.Range("G11").Value = 2.484984835592229
$ws.Range("H11").Value = 1.814290191105783
$ws.Range("L11").Value = 0.2864624392288988
$ws.Range("B12").Value = 2.441884499344326
$ws.Range("C12").Value = 0.2437212372405497
$ws.Range("D12").Value = 0.1153050421101085
$ws.Range("E12").Value = 0.07113679257140326
$ws.Range("G12").Value = 2.509760541113735
$ws.Range("H12").Value = 1.824407362856903
$ws.Range("L12").Value = 0.2896062891172306
$ws.Range("B13").Value = 2.433754677070851
$ws.Range("C13").Value = 0.2425367261945723
$ws.Range("D13").Value = 0.1147981805730041
$ws.Range("E13").Value = 0.07113327087970589
$ws.Range("G13").Value = 2.504414988250403
$ws.Range("H13").Value = 1.822222125259941
$ws.Range("L13").Value = 0.2889280780375714
$ws.Range("B14").Value = 2.40727538962966
$ws.Range("C14").Value = 0.2386727362544434
$ws.Range("D14").Value = 0.1131465745061888
$ws.Range("E14").Value = 0.07112231618891052
$ws.Range("G14").Value = 2.487018800413523
$ws.Range("H14").Value = 1.815119690108247
$ws.Range("L14").Value = 0.2867205782657152
$ws.Range("B15").Value = 2.391082021968259
$ws.Range("C15").Value = 0.236305124441941
$ws.Range("D15").Value = 0.1121359739476162
$ws.Range("E15").Value = 0.0711160134022375
$ws.Range("G15").Value = 2.476391349845358
$ws.Range("H15").Value = 1.810787727709851
$ws.Range("L15").Value = 0.2853717160594869
$ws.Range("B16").Value = 2.298680848942752
$ws.Range("C16").Value = 0.2227247034451807
$ws.Range("D16").Value = 0.1063606395876633
$ws.Range("E16").Value = 0.07108612793929581
$ws.Range("G16").Value = 2.415919010780414
$ws.Range("H16").Value = 1.786244902932651
$ws.Range("L16").Value = 0.2776922564144826
$ws.Range("B17").Value = 2.242347964615931
$ws.Range("C17").Value = 0.2143820640648357
$ws.Range("D17").Value = 0.1028317556032192
$ws.Range("E17").Value = 0.07107334627668749
$ws.Range("G17").Value = 2.379200778414884
$ws.Range("H17").Value = 1.77143869091779
$ws.Range("L17").Value = 0.2730258668335779
$ws.Range("B18").Value = 2.210075627995252
$ws.Range("C18").Value = 0.2095792768520823
$ws.Range("D18").Value = 0.1008071410741422
$ws.Range("E18").Value = 0.07106803105757109
$ws.Range("G18").Value = 2.358219550236811
$ws.Range("H18").Value = 1.763013826044158
$ws.Range("L18").Value = 0.2703582372120792
$ws.Range("B19").Value = 2.19917084317774
$ws.Range("C19").Value = 0.2079523841171351
$ws.Range("D19").Value = 0.1001225113808317
$ws.Range("E19").Value = 0.07106658080870076
$ws.Range("G19").Value = 2.351139261054669
$ws.Range("H19").Value = 1.760176931408239
$ws.Range("L19").Value = 0.2694578253225473
$ws.Range("B20").Value = 2.248331347625708
$ws.Range("C20").Value = 0.2152705973402362
$ws.Range("D20").Value = 0.1032068811888109
$ws.Range("E20").Value = 0.07107449605303451
$ws.Range("G20").Value = 2.383095170421598
$ws.Range("H20").Value = 1.773005377114373
$ws.Range("L20").Value = 0.273520917654011
$ws.Range("B21").Value = 2.415047584618321
$ws.Range("C21").Value = 0.2398078512528343
$ws.Range("D21").Value = 0.1136314707883628
$ws.Range("E21").Value = 0.0711254489487434
$ws.Range("G21").Value = 2.492122594431692
$ws.Range("H21").Value = 1.817201991991396
$ws.Range("L21").Value = 0.2873682873182588
$ws.Range("B22").Value = 2.525166584621275
$ws.Range("C22").Value = 0.2558080045952522
$ws.Range("D22").Value = 0.1204915994177611
$ws.Range("E22").Value = 0.07117696780766103
$ws.Range("G22").Value = 2.564637541667594
$ws.Range("H22").Value = 1.846912716816234
$ws.Range("L22").Value = 0.2965656194040776
$ws.Range("B23").Value = 2.466287193287712
$ws.Range("C23").Value = 0.2472716437804365
$ws.Range("D23").Value = 0.1168258363912713
$ws.Range("E23").Value = 0.07114780016237887
$ws.Range("G23").Value = 2.525818307354655
$ws.Range("H23").Value = 1.830979366023826
$ws.Range("L23").Value = 0.2916432736241887
$ws.Range("B24").Value = 2.245625904634551
$ws.Range("C24").Value = 0.2148689116937703
$ws.Range("D24").Value = 0.1030372738842686
$ws.Range("E24").Value = 0.0710739699074221
$ws.Range("G24").Value = 2.381334115642289
$ws.Range("H24").Value = 1.772296806122029
$ws.Range("L24").Value = 0.273297058038878
$ws.Range("B25").Value = 2.012525189863993
$ws.Range("C25").Value = 0.1797508009542526
$ws.Range("D25").Value = 0.08835833191416498
$ws.Range("E25").Value = 0.07107213990658323
$ws.Range("G25").Value = 2.23074782037142
$ws.Range("H25").Value = 1.712480115276492
$ws.Range("L25").Value = 0.2541319492743526
